# Sprint 3 acceptance testing (04/13) - fill in Pass/Fail + tester comments for Sprint 3
# column on the "Test Plan" worksheet for the first batch of user stories (rows 2-48).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Plan")
$ws.Activate()

$ws.Range("G2").Value = 'Pass'
$ws.Range("H2").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G3").Value = 'Pass'
$ws.Range("H3").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G4").Value = 'Pass'
$ws.Range("H4").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G5").Value = 'Pass'
$ws.Range("H5").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G6").Value = 'Pass'
$ws.Range("H6").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G7").Value = 'Pass'
$ws.Range("H7").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G8").Value = 'Pass'
$ws.Range("H8").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G9").Value = 'Pass'
$ws.Range("H9").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G10").Value = 'Pass'
$ws.Range("H10").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G11").Value = 'Pass'
$ws.Range("H11").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G12").Value = 'Pass'
$ws.Range("H12").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G13").Value = 'Pass'
$ws.Range("H13").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G14").Value = 'Pass'
$ws.Range("H14").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G15").Value = 'Pass'
$ws.Range("H15").Value = 'MC; 04/13; Further implementation completed and operational.'
$ws.Range("G16").Value = 'Pass'
$ws.Range("H16").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G17").Value = 'Pass'
$ws.Range("H17").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G18").Value = 'Pass'
$ws.Range("H18").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G19").Value = 'Pass'
$ws.Range("H19").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G20").Value = 'Pass'
$ws.Range("H20").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G21").Value = 'Pass'
$ws.Range("H21").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G22").Value = 'Pass'
$ws.Range("H22").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G23").Value = 'Pass'
$ws.Range("H23").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G24").Value = 'Fail'
$ws.Range("H24").Value = 'MC; 04/13; If a piece moves into a position in front of another jumpable piece, the "required jump" message is returned, and a second move must be made to jump the piece.'
$ws.Range("G25").Value = 'Fail'
$ws.Range("H25").Value = 'MC; 04/13; Required jumps are determined based on the end position of the piece, and so the above error is encountered. '
$ws.Range("G26").Value = 'Pass'
$ws.Range("H26").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G27").Value = 'Fail'
$ws.Range("H27").Value = 'MC; 04/13; A kinged piece may be moved on the same turn it was crowned if the second move is also valid. '
$ws.Range("G28").Value = 'Pass'
$ws.Range("H28").Value = 'MC; 04/13; Works as expected. '
$ws.Range("G29").Value = 'Pass'
$ws.Range("H29").Value = 'MC; 04/13; Still require connection to PlayerResignRoute to resign player upon signing out.'
$ws.Range("G30").Value = 'Fail'
$ws.Range("H30").Value = 'MC; 04/13; Implementation incomplete, player is not resigned from active game when signing out.'
$ws.Range("G31").Value = 'Fail'
$ws.Range("H31").Value = 'MC; 04/13; Implementation incomplete'
$ws.Range("G32").Value = 'Fail'
$ws.Range("H32").Value = 'MC; 04/13; Implementation incomplete'
$ws.Range("G33").Value = 'Fail'
$ws.Range("H33").Value = 'MC; 04/13; Implementation incomplete'
$ws.Range("G34").Value = 'Fail'
$ws.Range("H34").Value = 'MC; 04/13; Implementation incomplete'
$ws.Range("G35").Value = 'Fail'
$ws.Range("H35").Value = 'MC; 04/13; Implementation incomplete'
$ws.Range("G36").Value = 'Fail'
$ws.Range("H36").Value = 'MC; 04/13; Implementation incomplete'
$ws.Range("G37").Value = 'Fail'
$ws.Range("H37").Value = 'MC; 04/13; Implementation incomplete'
$ws.Range("G38").Value = 'Fail'
$ws.Range("H38").Value = 'MC; 04/13; Implementation incomplete'
$ws.Range("G39").Value = 'Fail'
$ws.Range("H39").Value = 'MC; 04/13; Implementation incomplete'
$ws.Range("G40").Value = 'Fail'
$ws.Range("H40").Value = 'MC; 04/13; Implementation incomplete'
$ws.Range("G41").Value = 'Fail'
$ws.Range("H41").Value = 'MC; 04/13; Implementation incomplete'
$ws.Range("G42").Value = 'Fail'
$ws.Range("H42").Value = 'MC; 04/13; Implementation incomplete'
$ws.Range("G43").Value = 'Fail'
$ws.Range("H43").Value = 'MC; 04/13; Implementation incomplete'
$ws.Range("G44").Value = 'Fail'
$ws.Range("H44").Value = 'MC; 04/13; Implementation incomplete'
$ws.Range("G45").Value = 'Fail'
$ws.Range("H45").Value = 'MC; 04/13; Implementation incomplete'
$ws.Range("G46").Value = 'Fail'
$ws.Range("H46").Value = 'MC; 04/13; Implementation incomplete'
$ws.Range("G47").Value = 'Fail'
$ws.Range("H47").Value = 'MC; 04/13; Implementation incomplete'
$ws.Range("G48").Value = 'Fail'
$ws.Range("H48").Value = 'MC; 04/13; Implementation incomplete'

# Leave the view scrolled/selected on the last cell touched during testing.
$ws.Range("G42").Select()
